# Update "想去人数" (F column) figures to the newly scraped counts.
$wb = $excel.ActiveWorkbook

# -- Sheet: 展览 (exhibitions) --
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 10712
$ws1.Range("F6").Value  = 1102
$ws1.Range("F9").Value  = 1897
$ws1.Range("F11").Value = 713
$ws1.Range("F12").Value = 213
$ws1.Range("F13").Value = 282
$ws1.Range("F14").Value = 253
$ws1.Range("F15").Value = 271
$ws1.Range("F16").Value = 967
$ws1.Range("F18").Value = 175
$ws1.Range("F22").Value = 166
$ws1.Range("F23").Value = 423
$ws1.Range("F24").Value = 177

# -- Sheet: 演出 (shows) --
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 614

# -- Sheet: 全部类型 (all types) --
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value  = 10713
$ws4.Range("F9").Value  = 1102
$ws4.Range("F13").Value = 1897
$ws4.Range("F15").Value = 713
$ws4.Range("F17").Value = 213
$ws4.Range("F18").Value = 282
$ws4.Range("F19").Value = 253
$ws4.Range("F20").Value = 271
$ws4.Range("F21").Value = 967
$ws4.Range("F23").Value = 614
$ws4.Range("F24").Value = 175
$ws4.Range("F29").Value = 166
$ws4.Range("F30").Value = 423
$ws4.Range("F31").Value = 177
